# Update handback/handoff timestamps as part of "Generate Report for Handback".
$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for 017a099c-...md (row 3, col G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-23 12:46:02"

# zh-cn sheet: "Correspond Handoff Datetime" (col H) and
# "Correspond Handback DateTime" (col K) for the 017a099c-... row (row 3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-23 12:45:55"
$wsZhCn.Range("K3").Value = "2016-08-23 12:46:29"

# de-de sheet: "Correspond Handoff Datetime" (col H) shares its value with the
# Overview sheet's "Latest HO Xliff Generate Date" for this row, so it also
# becomes 2016-08-23 12:46:02. "Correspond Handback DateTime" (col K) becomes
# 2016-08-23 12:46:36.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-23 12:46:02"
$wsDeDe.Range("K3").Value = "2016-08-23 12:46:36"
